$d = $word.ActiveDocument

$replacements = @(
    @{old="380÷7="; new="561÷4="},
    @{old="538÷9="; new="900÷8="},
    @{old="179÷6="; new="798÷8="},
    @{old="866÷8="; new="852÷6="},
    @{old="969÷7="; new="543÷2="},
    @{old="840÷6="; new="213÷5="},
    @{old="315÷8="; new="277÷7="},
    @{old="160÷7="; new="782÷4="},
    @{old="221÷8="; new="355÷9="},
    @{old="953÷9="; new="769÷4="},
    @{old="180÷5="; new="945÷7="},
    @{old="179÷9="; new="377÷4="},
    @{old="306÷7="; new="146÷6="},
    @{old="597÷5="; new="424÷6="},
    @{old="846÷4="; new="993÷6="},
    @{old="359÷5="; new="765÷6="},
    @{old="581÷4="; new="312÷5="},
    @{old="840÷8="; new="849÷4="},
    @{old="528÷4="; new="140÷3="},
    @{old="272÷7="; new="588÷2="},
    @{old="141÷6="; new="986÷5="},
    @{old="971÷6="; new="129÷8="},
    @{old="582÷4="; new="553÷7="},
    @{old="666÷4="; new="277÷3="},
    @{old="621÷3="; new="434÷2="}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
